# Add a "PRESUPUESTO" budget column (G) to the "VENTA MENSUAL" sheet, matching
# the layout/formatting of the existing month columns (A..F).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# New column width for column G. Excel stores column widths internally using a
# pixel-rounded "stored width" that is a bit wider than the character count you
# assign through ColumnWidth; 16.2 is the character width that round-trips to a
# stored width of exactly 17 (matching the target column width).
$ws.Columns.Item(7).ColumnWidth = 16.2

# Header cell G1: copy the header formatting (bold, centered, bordered) from F1,
# then set the label.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G1").Value = "PRESUPUESTO"

# Data rows G2:G29: copy the numeric/currency formatting from the F column data
# cells, then fill with 0 (matching the rest of the sheet's placeholder values).
$ws.Range("F2").Copy()
$ws.Range("G2:G29").PasteSpecial(-4122)  # xlPasteFormats
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 7).Value = 0
}

# Totals row G30: copy the totals-row formatting from F30, then set to 0.
$ws.Range("F30").Copy()
$ws.Range("G30").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G30").Value = 0
